# "ticker appears with results" - append a new results row (row 3) to the
# single data sheet, mirroring the structure of the existing row 2.
#
# Columns (row 1 headers): Date, ScoreFinal, Verdict, totalSentiment,
# wordCount, sentenceCount, posWordPercentage, negWordPercentage,
# posPhrasePercentage, negPhrasePercentage, ElapsedMs, posWordCount,
# negWordCount, positivePhraseCount, negativePhraseCount, Method, RSI,
# PEG, 200Moving%, 50Moving%, PriceBook, Dividend, Bollinger

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 3

$ws.Cells.Item($row, 1).Value  = (Get-Date -Year 2016 -Month 9 -Day 19 -Hour 21 -Minute 12 -Second 15)   # Date
$ws.Cells.Item($row, 2).Value  = -3        # ScoreFinal
$ws.Cells.Item($row, 3).Value  = "Buy"     # Verdict (new shared string)
$ws.Cells.Item($row, 4).Value  = 0         # totalSentiment
$ws.Cells.Item($row, 5).Value  = 0         # wordCount
$ws.Cells.Item($row, 6).Value  = 0         # sentenceCount
$ws.Cells.Item($row, 7).Value  = 0         # posWordPercentage
$ws.Cells.Item($row, 8).Value  = 0         # negWordPercentage
$ws.Cells.Item($row, 9).Value  = 0         # posPhrasePercentage
$ws.Cells.Item($row, 10).Value = 0         # negPhrasePercentage
$ws.Cells.Item($row, 11).Value = 0         # ElapsedMs
$ws.Cells.Item($row, 12).Value = 0         # posWordCount
$ws.Cells.Item($row, 13).Value = 0         # negWordCount
$ws.Cells.Item($row, 14).Value = 0         # positivePhraseCount
$ws.Cells.Item($row, 15).Value = 0         # negativePhraseCount
$ws.Cells.Item($row, 16).Value = "Random"  # Method
$ws.Cells.Item($row, 17).Value = 0         # RSI
$ws.Cells.Item($row, 18).Value = 0.87      # PEG
$ws.Cells.Item($row, 19).Value = 0.0351    # 200Moving%  (percentage format)
$ws.Cells.Item($row, 19).NumberFormat = "0.00%"
$ws.Cells.Item($row, 20).Value = -2.08     # 50Moving%
$ws.Cells.Item($row, 21).Value = 15.16     # PriceBook
$ws.Cells.Item($row, 22).Value = "N/A"     # Dividend
$ws.Cells.Item($row, 23).Value = 0         # Bollinger
